# Apply the "Added a design slide to final presentation" edit.
#
# Summary of the change:
#   1. Slide 3 ("VeriHandy Design") gains a new "Javascript" bullet after "CSS3".
#   2. A brand new slide ("Design" / "Two Types of Users") is inserted right
#      before the existing "Design" / "Six MySQL Tables" slide (position 5),
#      pushing the remaining slides (Six MySQL Tables, Today's Demo,
#      VeriHandy Demo) down by one.
#   3. The "Today's Demo" slide gets a couple of cosmetic run merges.

$p = $ppt.ActivePresentation

# --- 1. Slide 3: "VeriHandy Design" tech-stack list -> add "Javascript" ----
$techSlide = $p.Slides.Item(3)
$techBody = $techSlide.Shapes.Item(2).TextFrame.TextRange
[void]$techBody.InsertAfter("`rJavascript")

# --- 2. "Today's Demo" slide: merge "User " + "Registration/Login" runs ---
# (cosmetic cleanup that PowerPoint performs automatically when the text is
#  re-entered; re-assigning the full text forces the run list to collapse)
$demoSlide = $p.Slides.Item(6)
$demoBody = $demoSlide.Shapes.Item(2).TextFrame.TextRange
$demoFullText = $demoBody.Text
$demoBody.Text = ""
$demoBody.Text = $demoFullText

# --- 3. Insert the new "Design" / "Two Types of Users" slide at index 5 ---
$newSlide = $p.Slides.Add(5, 2)

$newSlide.Shapes.Item(1).TextFrame.TextRange.Text = "Design"

$bodyLines = @(
    @(0, "Two Types of Users"),
    @(1, "Normal User"),
    @(2, "Can be both a worker and a customer"),
    @(2, "Workers are those who complete jobs the customer users submit"),
    @(1, "Administrator User"),
    @(2, "Is able to view all users, jobs, and user messages"),
    @(2, "Also has normal user functionality")
)

$newBody = $newSlide.Shapes.Item(2).TextFrame.TextRange
$newBody.Text = ($bodyLines | ForEach-Object { $_[1] }) -join "`r"

$charPos = 1
foreach ($line in $bodyLines) {
    $level = $line[0]
    $text = $line[1]
    $len = $text.Length
    if ($level -gt 0) {
        $run = $newBody.Characters($charPos, $len)
        $run.IndentLevel = $level + 1
    }
    $charPos = $charPos + $len + 1
}
